$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 286, shifting existing rows 286-349 down to 287-350
$ws.Rows.Item(286).Insert()

# Fill in the new row 286 with the fresh data record
$ws.Cells.Item(286, 1).Value = 4
$ws.Cells.Item(286, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(286, 3).Value = "Los Lagos"
$ws.Cells.Item(286, 4).Value = 44889
$ws.Cells.Item(286, 5).Value = 10
$ws.Cells.Item(286, 6).Value = 100112017
$ws.Cells.Item(286, 7).Value = "Apio"
$ws.Cells.Item(286, 8).Value = "Americana (o)"
$ws.Cells.Item(286, 9).Value = "Primera"
$ws.Cells.Item(286, 10).Value = 50
$ws.Cells.Item(286, 11).Value = 15000
$ws.Cells.Item(286, 12).Value = 15000
$ws.Cells.Item(286, 13).Value = 15000
$ws.Cells.Item(286, 14).Value = "`$/docena de matas"
$ws.Cells.Item(286, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(286, 16).Value = 2500
$ws.Cells.Item(286, 17).Value = 6
$ws.Cells.Item(286, 18).Value = "Hortaliza"
